$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete review rows (old rows 2-4). The two rows that follow
# (old rows 5 and 6) shift up to become the new rows 2 and 3.
$ws.Rows("2:4").Delete()

# The removed rows carried hyperlinks on column D; drop them now that the
# corresponding rows are gone.
$ws.Hyperlinks.Delete()
